# "added encoded password utility"
# The shared-string value that used to hold an (encoded) password is
# replaced with the plain config-key reference, and column B is widened
# so the longer text keeps fitting/displaying cleanly ("bestFit").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 held the old encoded password value ("Q2ZnQDIwMjEk") - swap it for
# the new config-reference string. Leading apostrophe keeps it entered
# as literal text (preserves the cell's existing "quote prefix" style)
# without ending up in the stored value itself.
$ws.Range("B2").Value = "'config:approverpassword"

# Column B needs to be a bit wider to comfortably fit the new, longer
# value (old bestFit width ~15.86 chars -> new ~23.86 chars).
$ws.Columns.Item(2).ColumnWidth = 23
